$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear CommitID values that were posted for first run (A2, A3), but keep formatting/style
$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()

# Update Whitelist value in row 2 from old host entry to new host entry
$ws.Range("E2").Value = "ip-10-123-10-250"

# Update the active selection to match saved view state
$ws.Range("B5").Select()

$wb.Save()
